# Append new time-tracking entries to the status report (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each entry: day serial (Excel 1900 date system), hours, comment text.
$newRows = @(
    @(40260, 2.5, "Group Meeting"),
    @(40260, 1,   "Weekly Meeting"),
    @(40261, 3.5, "Physical interface board design"),
    @(40262, 2,   "QNX - Encoder Test"),
    @(40263, 0.5, "Skype Meeting"),
    @(40263, 6,   "QNX - Encoder Test"),
    @(40264, 8,   "QNX - Encoder Test")
)

$startRow = 65
$lastExistingRow = 64

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]

    # Carry the existing date-column formatting down onto the new rows
    # (same style used by every other data row), instead of inventing a
    # brand-new number format.
    $ws.Range("A" + $lastExistingRow).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# Keep the view in sync with the newly-added rows (matches Excel's own
# behaviour of scrolling/selecting near the last edited cell).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 43
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A70").Select()
